$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.782.48'
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('D3').Value = '2.468.74'
$ws.Range('E3').Value = '  -0.58%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '316.46'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.36%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '93.00'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.550'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +1.03%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +3.81%  '
$ws.Range('E10').Value = '  +0.94%  '
$ws.Range('E11').Value = '  +8.88%  '
$ws.Range('E12').Value = '  +0.16%  '
$ws.Range('D13').Value = '2.849.52'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('E14').Value = '  +0.88%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.74'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.20%  '
$ws.Range('D16').Value = '2.472.80'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.781'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +3.61%  '
$ws.Range('D18').Value = '41.751.16'
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.51'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +3.29%  '
$ws.Range('E20').Value = '  +3.33%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.56'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.87%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '71.15'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.96%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '239.60'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.88%  '
$ws.Range('E24').Value = '  +0.50%  '
$ws.Range('E25').Value = '  +1.32%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '24.84'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.78'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.72%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '36.13'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '156.12'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.88%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.53'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.42%  '
$ws.Range('E33').Value = '  +0.13%  '
$ws.Range('E34').Value = '  +1.78%  '
$ws.Range('E35').Value = '  +1.26%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '17.61'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.70%  '
$ws.Range('E37').Value = '  -1.63%  '
$ws.Range('E38').Value = '  +1.43%  '
$ws.Range('E39').Value = '  -0.43%  '
$ws.Range('E40').Value = '  -0.89%  '
$ws.Range('E41').Value = '  -2.76%  '
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('D43').Value = '1.973.44'
$ws.Range('E43').Value = '  +1.11%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '18.99'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.75%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0284'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('E47').Value = '  +2.14%  '
$ws.Range('D48').Value = '2.702.99'
$ws.Range('E48').Value = '  -0.82%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '96.93'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.01%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '67.30'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.29%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '73.01'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.20%  '
